$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46016
$ws.Range("D8").Value = 152.16
$ws.Range("E8").Value = 151.68
$ws.Range("F8").Value = 161.68
$ws.Range("G8").Value = 151.79
$ws.Range("A9").Value = 46016
$ws.Range("D9").Value = 152.16
$ws.Range("E9").Value = 151.68
$ws.Range("F9").Value = 161.68
$ws.Range("G9").Value = 151.79
$ws.Range("A10").Value = 46016
$ws.Range("D10").Value = 154.87
$ws.Range("E10").Value = 153.79
$ws.Range("F10").Value = 163.79
$ws.Range("G10").Value = 154.31
$ws.Range("A11").Value = 46015
$ws.Range("D11").Value = 152.34
$ws.Range("E11").Value = 152.05000000000001
$ws.Range("F11").Value = 162.05000000000001
$ws.Range("G11").Value = 152.16999999999999
$ws.Range("A12").Value = 46015
$ws.Range("D12").Value = 152.34
$ws.Range("E12").Value = 152.05000000000001
$ws.Range("F12").Value = 162.05000000000001
$ws.Range("G12").Value = 152.16999999999999
$ws.Range("A13").Value = 46015
$ws.Range("D13").Value = 155.07
$ws.Range("E13").Value = 154.18
$ws.Range("F13").Value = 164.18
$ws.Range("G13").Value = 154.69999999999999
$ws.Range("A17").Value = 46016
$ws.Range("D17").Value = 159.04
$ws.Range("E17").Value = 159.72999999999999
$ws.Range("F17").Value = 169.73
$ws.Range("A18").Value = 46015
$ws.Range("D18").Value = 159.21
$ws.Range("E18").Value = 160.13999999999999
$ws.Range("F18").Value = 170.14
$ws.Range("A22").Value = 46016
$ws.Range("D22").Value = 153.74
$ws.Range("E22").Value = 153.05000000000001
$ws.Range("F22").Value = 162.65
$ws.Range("G22").Value = 154.21
$ws.Range("A23").Value = 46016
$ws.Range("D23").Value = 159.33000000000001
$ws.Range("E23").Value = 158.66
$ws.Range("F23").Value = 168.66
$ws.Range("G23").Value = "N/A"
$ws.Range("A24").Value = 46016
$ws.Range("D24").Value = 159.12
$ws.Range("E24").Value = 159.16
$ws.Range("F24").Value = 169.16
$ws.Range("G24").Value = "N/A"
$ws.Range("A25").Value = 46016
$ws.Range("D25").Value = 159.62
$ws.Range("E25").Value = 158.9
$ws.Range("F25").Value = 168.9
$ws.Range("G25").Value = 158.66999999999999
$ws.Range("A26").Value = 46016
$ws.Range("D26").Value = 158.56
$ws.Range("E26").Value = 160.18
$ws.Range("F26").Value = 170.18
$ws.Range("G26").Value = "N/A"
$ws.Range("A27").Value = 46015
$ws.Range("D27").Value = 153.93
$ws.Range("E27").Value = 153.43
$ws.Range("F27").Value = 163.03
$ws.Range("G27").Value = 154.59
$ws.Range("A28").Value = 46015
$ws.Range("D28").Value = 159.53
$ws.Range("E28").Value = 159.05000000000001
$ws.Range("F28").Value = 169.05
$ws.Range("G28").Value = "N/A"
$ws.Range("A29").Value = 46015
$ws.Range("D29").Value = 159.31
$ws.Range("E29").Value = 159.56
$ws.Range("F29").Value = 169.56
$ws.Range("G29").Value = "N/A"
$ws.Range("A30").Value = 46015
$ws.Range("D30").Value = 159.81
$ws.Range("E30").Value = 159.30000000000001
$ws.Range("F30").Value = 169.3
$ws.Range("G30").Value = 159.07
$ws.Range("A31").Value = 46015
$ws.Range("D31").Value = 158.74
$ws.Range("E31").Value = 160.58000000000001
$ws.Range("F31").Value = 170.58
$ws.Range("G31").Value = "N/A"
$ws.Range("A35").Value = 46016
$ws.Range("D35").Value = 152.38999999999999
$ws.Range("E35").Value = 152.24
$ws.Range("F35").Value = 161.24
$ws.Range("A36").Value = 46015
$ws.Range("D36").Value = 152.58000000000001
$ws.Range("E36").Value = 152.63
$ws.Range("F36").Value = 161.63
$ws.Range("A40").Value = 46016
$ws.Range("D40").Value = 160.15
$ws.Range("E40").Value = 160.18
$ws.Range("F40").Value = 170.18
$ws.Range("A41").Value = 46016
$ws.Range("D41").Value = 159.87
$ws.Range("E41").Value = 160.6
$ws.Range("F41").Value = 170.6
$ws.Range("A42").Value = 46015
$ws.Range("D42").Value = 160.37
$ws.Range("E42").Value = 160.63
$ws.Range("F42").Value = 170.63
$ws.Range("A43").Value = 46015
$ws.Range("D43").Value = 160.09
$ws.Range("E43").Value = 161.05000000000001
$ws.Range("F43").Value = 171.05
$ws.Range("A47").Value = 46016
$ws.Range("D47").Value = 153.63999999999999
$ws.Range("E47").Value = 154.47
$ws.Range("F47").Value = 164.47
$ws.Range("A48").Value = 46016
$ws.Range("D48").Value = 153.47
$ws.Range("E48").Value = 154.57
$ws.Range("F48").Value = 164.57
$ws.Range("A49").Value = 46015
$ws.Range("D49").Value = 154.16999999999999
$ws.Range("E49").Value = 155.25
$ws.Range("F49").Value = 165.25
$ws.Range("A50").Value = 46015
$ws.Range("D50").Value = 153.99
$ws.Range("E50").Value = 155.35
$ws.Range("F50").Value = 165.35
$ws.Range("A54").Value = 46016
$ws.Range("D54").Value = 169.79
$ws.Range("E54").Value = 169.5
$ws.Range("F54").Value = 179.5
$ws.Range("A55").Value = 46016
$ws.Range("D55").Value = 157.93
$ws.Range("E55").Value = 164.26
$ws.Range("F55").Value = 174.26
$ws.Range("A56").Value = 46016
$ws.Range("D56").Value = 159.63
$ws.Range("E56").Value = "N/A"
$ws.Range("F56").Value = "N/A"
$ws.Range("A57").Value = 46016
$ws.Range("D57").Value = 158.81
$ws.Range("E57").Value = 158.54
$ws.Range("F57").Value = "N/A"
$ws.Range("A58").Value = 46016
$ws.Range("D58").Value = 154.71
$ws.Range("E58").Value = 154.58000000000001
$ws.Range("F58").Value = 164.58
$ws.Range("A59").Value = 46016
$ws.Range("D59").Value = 161.97
$ws.Range("E59").Value = 166.97
$ws.Range("F59").Value = "N/A"
$ws.Range("A60").Value = 46015
$ws.Range("D60").Value = 169.98
$ws.Range("E60").Value = 169.99
$ws.Range("F60").Value = 179.99
$ws.Range("A61").Value = 46015
$ws.Range("D61").Value = 158.13
$ws.Range("E61").Value = 164.62
$ws.Range("F61").Value = 174.62
$ws.Range("A62").Value = 46015
$ws.Range("D62").Value = 159.83000000000001
$ws.Range("E62").Value = "N/A"
$ws.Range("F62").Value = "N/A"
$ws.Range("A63").Value = 46015
$ws.Range("D63").Value = 158.97999999999999
$ws.Range("E63").Value = 158.9
$ws.Range("F63").Value = "N/A"
$ws.Range("A64").Value = 46015
$ws.Range("D64").Value = 154.88
$ws.Range("E64").Value = 154.94
$ws.Range("F64").Value = 164.94
$ws.Range("A65").Value = 46015
$ws.Range("D65").Value = 162.13999999999999
$ws.Range("E65").Value = 167.41
$ws.Range("F65").Value = "N/A"
